$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: give the row a custom (taller) height ---
$ws.Rows.Item(11).RowHeight = 15.65

# --- Column width tweaks (C, E, F) ---
$ws.Columns.Item(3).ColumnWidth = 30.2
$ws.Columns.Item(5).ColumnWidth = 20.2
$ws.Columns.Item(6).ColumnWidth = 21

# --- Section 3 input columns (replaces old "Section 3" block) ---
$ws.Range("C25").Value = "Section 3 input columns"
$ws.Range("C25").Font.Bold = $true
$ws.Range("C25").Font.Underline = $true

$ws.Range("C26").Value = "Operation"
$ws.Range("C26").Font.Bold = $true
$ws.Range("D26").Value = "Signal 1"
$ws.Range("D26").Font.Bold = $true
$ws.Range("E26").Value = "Signal 2"
$ws.Range("E26").Font.Bold = $true
$ws.Range("F26").Value = "Signal 3"
$ws.Range("F26").Font.Bold = $true

$ws.Range("C27").Value = "(none,  and, or)"
$ws.Range("D27").Value = "mandatory"
$ws.Range("E27").Value = "Only if (and, or)"
$ws.Range("F27").Value = "Only if (and, or)"

$ws.Range("B28").Value = "Col 0"
$ws.Range("C28").Value = "none"
$ws.Range("C28").NumberFormat = """TRUE"";""TRUE"";""FALSE"""
$ws.Range("D28").Value = "variable crossPercent 3"

$ws.Range("B29").Value = "Col 1"
$ws.Range("C29").Value = "or"
$ws.Range("D29").Value = "crossAbove 10 50"
$ws.Range("E29").Value = "crossAbove 10 100"
$ws.Range("F29").Value = "crossAbove 10 200"

$ws.Range("B30").Value = "Col 2"
$ws.Range("C30").Value = "none"
$ws.Range("D30").Value = "topLine 10"

$ws.Range("B31").Value = "Col 3"
$ws.Range("C31").Value = "and"
$ws.Range("D31").Value = "priceAbove 10 100"
$ws.Range("E31").Value = "priceAbove 10 200"

$ws.Range("B32").Value = "Col 4"
$ws.Range("C32").Value = "none"
$ws.Range("D32").Value = "priceAbove 50 30"

$ws.Range("B33").Value = "Col 5"
$ws.Range("B34").Value = "Col 6*"
$ws.Range("B35").Value = "Col 7*"
$ws.Range("B36").Value = "Col 8*"
$ws.Range("B37").Value = "Col 9*"

$ws.Range("B39").Value = "* not mandatory"

# --- Section 3 input test ---
$ws.Range("C42").Value = "Section 3 input test"
$ws.Range("C42").Font.Bold = $true
$ws.Range("C42").Font.Underline = $true

$ws.Range("C43").Value = "Operation"
$ws.Range("C43").Font.Bold = $true
$ws.Range("D43").Value = "Column"
$ws.Range("D43").Font.Bold = $true
$ws.Range("E43").Value = "Column"
$ws.Range("E43").Font.Bold = $true
$ws.Range("F43").Value = "Column"
$ws.Range("F43").Font.Bold = $true

$ws.Range("A45").Value = "Test 1"
$ws.Range("B45").Value = "Part 1"
$ws.Range("B46").Value = "Part 2"
$ws.Range("B47").Value = "Part 3"
$ws.Range("B48").Value = "Part 4"
$ws.Range("B49").Value = "Part 5"
$ws.Range("B50").Value = "Part 6"
$ws.Range("B51").Value = "Part 7"
$ws.Range("B52").Value = "Part 8"

$ws.Range("A54").Value = "Test 2"
$ws.Range("B54").Value = "Part 1"
$ws.Range("B55").Value = "Part 2"
$ws.Range("B56").Value = "Part 3"
$ws.Range("B57").Value = "Part 4"

# --- Keep the on-screen selection / scroll position in sync with the edit ---
$ws.Range("C54").Select()
